$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "Lienhart._ First" -> "Lienhart. First"
# Remove the stray underscore right after the period, splitting the run into
# two runs ("." and " First, a classifier...") the way Word does when you
# click right after the "." and press Delete.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Lienhart._ First")
if (-not $found) { throw "Hunk1: anchor text not found" }
$underscorePos = $rng.Start + 9   # "Lienhart." is 9 characters; "_" follows immediately
$splitPoint = $d.Range($underscorePos, $underscorePos)
$d.Bookmarks.Add("TmpSplit1", $splitPoint)
$d.Range($underscorePos, $underscorePos + 1).Delete()
$d.Bookmarks("TmpSplit1").Delete()

# ---------------------------------------------------------------------------
# Hunk 2: the "_GoBack" bookmark (Word's "last edit location" marker) moves
# from the end of the paragraph into the middle of the word "few" ("fe|w"),
# splitting that run in two. The visible text is unchanged.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("a few hundreds")
if (-not $found2) { throw "Hunk2: anchor text not found" }
$splitPos = $rng2.Start + 4   # "a fe" is 4 characters
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Hunk 3: "lntel@ lntegrated" -> "lntel lntegrated" (drop the stray "@").
# Delete the "@" using a "backspace from just after it" pattern so the
# surviving run keeps its original rsidRPr/attributes.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found3 = $rng3.Find.Execute("lntel@ lntegrated")
if (-not $found3) { throw "Hunk3: anchor text not found" }
$afterAt3 = $rng3.Start + 6   # "lntel@" is 6 characters
$splitPoint3 = $d.Range($afterAt3, $afterAt3)
$d.Bookmarks.Add("TmpSplit3", $splitPoint3)
$d.Range($afterAt3 - 1, $afterAt3).Delete()
$d.Bookmarks("TmpSplit3").Delete()

# ---------------------------------------------------------------------------
# Hunk 4: "lntel@ IPP) is an extensive" -> "lntel IPP) is an extensive"
# (drop the other stray "@"), same backspace pattern.
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$found4 = $rng4.Find.Execute("lntel@ IPP) is an extensive")
if (-not $found4) { throw "Hunk4: anchor text not found" }
$afterAt4 = $rng4.Start + 6   # "lntel@" is 6 characters
$splitPoint4 = $d.Range($afterAt4, $afterAt4)
$d.Bookmarks.Add("TmpSplit4", $splitPoint4)
$d.Range($afterAt4 - 1, $afterAt4).Delete()
$d.Bookmarks("TmpSplit4").Delete()

# ---------------------------------------------------------------------------
# Hunk 5: the cached PAGE field result in the header goes from "29" to "30".
# Edit the field result text in place (per-character, which is the one
# approach that reliably targets the field-result run in this runtime).
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
foreach ($f in $hdr.Range.Fields) {
    if ($f.Result.Text -eq "29") {
        $f.Result.Characters.Item(1).Text = "3"
        $f.Result.Characters.Item(2).Text = "0"
    }
}

Write-Output "done"
